$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.021.63"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "2.301.32"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'302.39"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "'99.11"
$ws.Range("E6").Value = "  +5.71%  "
$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.506"
$ws.Range("E9").Value = "  +2.95%  "
$ws.Range("D10").Value = "'34.34"
$ws.Range("E10").Value = "  +4.24%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "'49.23"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("E13").Value = "  +4.30%  "
$ws.Range("D14").Value = "'18.00"
$ws.Range("E14").Value = "  +18.00%  "
$ws.Range("D15").Value = "'6.79"
$ws.Range("D16").Value = "2.661.62"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "2.285.27"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "'0.809"
$ws.Range("E18").Value = "  +4.78%  "
$ws.Range("D19").Value = "42.918.59"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "'12.29"
$ws.Range("E20").Value = "  +7.93%  "
$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").Value = "'67.86"
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("D24").Value = "'236.53"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = "  +13.61%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.46"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").Value = "'24.78"
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("D29").Value = "'168.20"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "'2.09"
$ws.Range("E30").Value = "  -9.26%  "
$ws.Range("D31").Value = "'33.79"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'9.15"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("E36").Value = "  +3.89%  "
$ws.Range("D37").Value = "'16.91"
$ws.Range("E37").Value = "  +6.06%  "
$ws.Range("D38").Value = "'0.0699"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +3.70%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.81"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.78"
$ws.Range("E41").Value = "  +4.36%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").Value = "1.998.85"
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("D45").Value = "'0.0285"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").Value = "'10.04"
$ws.Range("E46").Value = "  +5.04%  "
$ws.Range("D47").Value = "'17.57"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "'2.86"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").Value = "'55.14"
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("D50").Value = "2.527.03"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("E51").Value = "  +2.63%  "
